$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data appended at row 48 (rows 45 and 47 are intentionally
# left blank, matching the existing sheet1 layout where row 45 is also
# blank between rows 44 and 46).
$ws.Range("A48").Value = "Final VM  and Retiered records NAME do not match to concepts/synonyms for multiple concepts"
$ws.Range("B48").Value = 3812991
$ws.Range("C48").Value = "E1CC4283-D134-A91B-E040-BB89AD43126C"
$ws.Range("D48").Value = "C107221:C25301"
$ws.Range("E48").Value = "Several Day"
$ws.Range("F48").Value = "Several days"
$ws.Range("G48").Value = "E141992A-68C7-79DB-E040-BB89AD43560A"
$ws.Range("H48").Value = 3811310
$ws.Range("I48").Value = "E141992A-68CA-79DB-E040-BB89AD43560A"
$ws.Range("J48").Value = "D9344734-8CAF-4378-E034-0003BA12F5E7 1 VM Alt Name"

$ws.Range("A48:J48").Style = "Normal"
$ws.Range("A48:J48").RowHeight = 58
$ws.Range("A48:J48").HorizontalAlignment = -4131
$ws.Range("A48:J48").VerticalAlignment = -4160
$ws.Range("A48:J48").WrapText = $true
$ws.Range("A48:J48").NumberFormat = "@"

$ws.Range("J51").Select()
$ws.Application.ActiveWindow.ScrollRow = 44
